# Applies the commit "adding averages and more checks":
#  - Header/title font changes: no longer a distinct bold-14pt font for titles;
#    both the title and the table-header fonts become bold, default size, white text.
#  - The "LAST UPDATE" date used to recompute "PERIOD TO EXPIRE" moves from
#    08-Sep-2025 to 16-Sep-2025 on the Training Dashboard, shifting every
#    "PERIOD TO EXPIRE" value down by 8 and updating "LAST UPDATE" text.
#  - Because of that date shift, the "LOTO (SOPs)" row (row 21) now falls
#    below the validity threshold, so it is flagged NOT VALID and highlighted
#    like the other invalid rows.
#  - The Exam Dashboard's COMMENTS column is widened and its comments text
#    is reworded from "OK" to "date is valid".

$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1. Font styling: titles (A1) and table header rows lose the old bold/14pt
#    font and instead share one bold, default-size, white font.
# ---------------------------------------------------------------------------
$white = 16777215

$trainingWs.Range("A1").Font.Color = $white
$trainingWs.Range("A1").Font.Size = 11

$examWs.Range("A1").Font.Color = $white
$examWs.Range("A1").Font.Size = 11

$trainingWs.Range("A2:K2").Font.Color = $white
$trainingWs.Range("A2:K2").Font.Size = 11

$examWs.Range("A2:G2").Font.Color = $white
$examWs.Range("A2:G2").Font.Size = 11

# ---------------------------------------------------------------------------
# 2. Training Dashboard: shift LAST UPDATE from 08-Sep-2025 to 16-Sep-2025,
#    and recompute PERIOD TO EXPIRE (down by 8 for every row).
# ---------------------------------------------------------------------------
$newLastUpdate = "16-Sep-2025"

$periodToExpire = @{
    3  = 675
    4  = 382
    5  = 360
    6  = 358
    7  = 405
    8  = 324
    9  = 406
    10 = 342
    11 = 344
    12 = 719
    13 = 446
    14 = 409
    15 = 408
    16 = 327
    17 = 386
    18 = 426
    19 = 427
    20 = 502
    21 = 15
    22 = -104
    23 = -190
    24 = 228
    25 = -45
    26 = 182
    27 = 200
    28 = 182
    29 = 213
    30 = 199
    31 = 228
    32 = 315
    33 = 313
    34 = 313
    35 = 313
    36 = 354
    37 = 354
}

# Keep the LAST UPDATE column stored as literal text (matching the rest of
# the sheet's date-like columns) instead of letting Excel auto-convert the
# string into a date serial number.
$trainingWs.Range("I3:I37").NumberFormat = "@"

foreach ($r in 3..37) {
    $trainingWs.Cells.Item($r, 8).Value = $periodToExpire[$r]
    $trainingWs.Cells.Item($r, 9).Value = $newLastUpdate
}

# ---------------------------------------------------------------------------
# 3. Row 21 (LOTO (SOPs)) now fails validity: update its STATUS text and
#    reformat the row like the other NOT VALID rows (row 22's formatting).
# ---------------------------------------------------------------------------
$trainingWs.Range("J21").Value = "NOT VALID"

$trainingWs.Range("A22:K22").Copy()
$trainingWs.Range("A21:K21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Exam Dashboard: widen COMMENTS column and reword the comments.
# ---------------------------------------------------------------------------
# ColumnWidth is expressed in characters-of-the-default-font and is offset
# from the stored OOXML column width by ~0.83; 14.17 yields a saved width of
# 15, matching the target column width.
$examWs.Columns.Item(5).ColumnWidth = 14.17

$examWs.Range("E3").Value = "date is valid"
$examWs.Range("E4").Value = "date is valid"
$examWs.Range("E5").Value = "date is valid"
